$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 129; this shifts rows 129-162 down to 130-163
# (Excel-style insert, copies formatting from the row above by default,
# which keeps the date numeric format on column D)
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new weekly record
$row = 129
$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value  = "Maule"
$ws.Cells.Item($row, 4).Value  = 44508
$ws.Cells.Item($row, 5).Value  = 7
$ws.Cells.Item($row, 6).Value  = 100112008
$ws.Cells.Item($row, 7).Value  = "Coliflor"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 5000
$ws.Cells.Item($row, 11).Value = 500
$ws.Cells.Item($row, 12).Value = 500
$ws.Cells.Item($row, 13).Value = 500
$ws.Cells.Item($row, 14).Value = "$/unidad"
$ws.Cells.Item($row, 15).Value = "Región del Maule"
$ws.Cells.Item($row, 16).Value = 500
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"

"Row inserted at 129 and populated"
